# Generate Report for Handoff
# - Priority for the four "low" priority files is bumped to "ht" (hot) on both the
#   zh-cn and de-de handoff status sheets.
# - The corresponding "Latest Handoff Datetime" timestamps are refreshed to reflect
#   the new handoff generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn: rows 4-7 -> Priority (E) low -> ht, Latest Handoff Datetime (H) refreshed
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("E5").Value = "ht"
$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("E7").Value = "ht"

$wsZhCn.Range("H4").Value = "2016-09-05 14:39:49"
$wsZhCn.Range("H5").Value = "2016-09-05 14:39:49"
$wsZhCn.Range("H6").Value = "2016-09-05 14:39:49"
$wsZhCn.Range("H7").Value = "2016-09-05 14:39:49"

# de-de: rows 4-7 -> Priority (E) low -> ht, Latest Handoff Datetime (H) refreshed
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("E7").Value = "ht"

$wsDeDe.Range("H4").Value = "2016-09-05 14:39:54"
$wsDeDe.Range("H5").Value = "2016-09-05 14:39:54"
$wsDeDe.Range("H6").Value = "2016-09-05 14:39:54"
$wsDeDe.Range("H7").Value = "2016-09-05 14:39:54"

# Overview: rows 4-7 -> Latest HO Xliff Generate Date (G) mirrors the de-de refresh
# (it tracks the most recent handoff generation across all target languages)
$wsOverview.Range("G4").Value = "2016-09-05 14:39:54"
$wsOverview.Range("G5").Value = "2016-09-05 14:39:54"
$wsOverview.Range("G6").Value = "2016-09-05 14:39:54"
$wsOverview.Range("G7").Value = "2016-09-05 14:39:54"
